$wb = $excel.ActiveWorkbook

# ===================================================================
# Nike sheet: replace the "LEVELUP / 60% off" promo rows with the new
# "personal growth / innovation / statement product" advantages.
# Old sheet had 6 data rows (2-7); new sheet has 5 data rows (2-6).
# ===================================================================
$nike = $wb.Worksheets.Item("Nike")

$nike.Cells.Item(2, 1).Value = 1
$nike.Cells.Item(2, 2).Value = 'Nike as a brand of personal growth'
$nike.Cells.Item(2, 3).Value = 'Brand empowerment'
$nike.Cells.Item(2, 4).Value = 'Nike continues to position itself as more than a product brand, emphasizing personal growth, individual choices, and empowerment. | This resonates strongly with consumers who value authenticity and personal journeys.'
$nike.Cells.Item(2, 5).ClearContents()
$nike.Cells.Item(2, 6).Value = 'Nike allowed me to retrace my growth over the past 10 years. Where I am and who I am today are the result of choices — and of never accepting compromises. I’m an artisan, not an influencer. #nikeshox #teamnike #ad'

$nike.Cells.Item(3, 1).Value = 2
$nike.Cells.Item(3, 2).Value = 'Revolutionary product design'
$nike.Cells.Item(3, 3).Value = 'Innovation'
$nike.Cells.Item(3, 4).Value = 'Nike''s focus on revolutionary technology and cutting-edge product design continues to resonate. | Ads emphasize how their products go beyond style to offer complete experiences.'
$nike.Cells.Item(3, 5).ClearContents()
$nike.Cells.Item(3, 6).Value = 'DNA rivoluzionario. Esperienza totale.'

$nike.Cells.Item(4, 1).Value = 2
$nike.Cells.Item(4, 2).Value = 'Revolutionary product design'
$nike.Cells.Item(4, 3).Value = 'Innovation'
$nike.Cells.Item(4, 4).Value = 'Nike''s focus on revolutionary technology and cutting-edge product design continues to resonate. | Ads emphasize how their products go beyond style to offer complete experiences.'
$nike.Cells.Item(4, 5).ClearContents()
$nike.Cells.Item(4, 6).Value = 'It’s not just about style. It’s the energy you carry when you put them on. #nikeshox #teamnike #ad'

$nike.Cells.Item(5, 1).Value = 3
$nike.Cells.Item(5, 2).Value = 'Nike Shox as a statement product'
$nike.Cells.Item(5, 3).Value = 'Statement product'
$nike.Cells.Item(5, 4).Value = 'Nike Shox is continuously pushed as a signature product with distinct branding, showcasing its power and innovation. | This taps into the desire for high-impact, noticeable products.'
$nike.Cells.Item(5, 5).ClearContents()
$nike.Cells.Item(5, 6).Value = 'Née pour déranger. L’expérience totale. (Born to disturb. The total experience.)'

$nike.Cells.Item(6, 1).Value = 3
$nike.Cells.Item(6, 2).Value = 'Nike Shox as a statement product'
$nike.Cells.Item(6, 3).Value = 'Statement product'
$nike.Cells.Item(6, 4).Value = 'Nike Shox is continuously pushed as a signature product with distinct branding, showcasing its power and innovation. | This taps into the desire for high-impact, noticeable products.'
$nike.Cells.Item(6, 5).ClearContents()
$nike.Cells.Item(6, 6).Value = 'La nouvelle Nike Shok Z, la paire qui ne te laisse pas indifférent ✨@nike #NikeShox #TeamNike Publicité 🎥 : @mehdiscovers (The new Nike Shok Z, the pair that doesn’t leave you indifferent ✨@nike #NikeShox #TeamNike Advertisement 🎥 : @mehdiscovers)'

$nike.Rows.Item(7).Delete()

# ===================================================================
# adidas sheet: drop the old "Extra 35% with promo code" advantage and
# add the new "Black Friday/adiClub" and "adidas Vibes Fragrance" rows.
# Old sheet had 4 data rows (2-5); new sheet has 5 data rows (2-6).
# ===================================================================
$adidas = $wb.Worksheets.Item("adidas")

$adidas.Cells.Item(2, 1).Value = 1
$adidas.Cells.Item(2, 2).Value = 'Iconic adidas style'
$adidas.Cells.Item(2, 3).Value = 'Brand/style positioning'
$adidas.Cells.Item(2, 4).Value = 'Several ads focus on taking your style up a notch with unmistakable adidas looks in different languages. | This stresses recognisable brand-led styling over price.'
$adidas.Cells.Item(2, 5).ClearContents()
$adidas.Cells.Item(2, 6).Value = 'Take it up a notch in unmistakable adidas style'

$adidas.Cells.Item(3, 1).Value = 1
$adidas.Cells.Item(3, 2).Value = 'Iconic adidas style'
$adidas.Cells.Item(3, 3).Value = 'Brand/style positioning'
$adidas.Cells.Item(3, 4).Value = 'Several ads focus on taking your style up a notch with unmistakable adidas looks in different languages. | This stresses recognisable brand-led styling over price.'
$adidas.Cells.Item(3, 5).ClearContents()
$adidas.Cells.Item(3, 6).Value = 'Schalt mit unverwechselbaren adidas Styles einen Gang höher (Switch up your style with unmistakable adidas styles)'

$adidas.Cells.Item(4, 1).Value = 2
$adidas.Cells.Item(4, 2).Value = 'Black Friday / adiClub Sign-up Incentives'
$adidas.Cells.Item(4, 3).Value = 'Promotional exclusivity'
$adidas.Cells.Item(4, 4).Value = 'Several ads offer an incentive to sign up for adiClub to participate in exclusive offers like winning a gift card. | This creates urgency around registration and exclusivity.'
$adidas.Cells.Item(4, 5).ClearContents()
$adidas.Cells.Item(4, 6).Value = 'In occasione del Black Friday, iscriviti ad adiClub e crea la tua wishlist per avere la possibilità di vincere una Gift Card da 250 €. (On the occasion of Black Friday, sign up for adiClub and create your wishlist for a chance to win a €250 Gift Card.)'

$adidas.Cells.Item(5, 1).Value = 2
$adidas.Cells.Item(5, 2).Value = 'Black Friday / adiClub Sign-up Incentives'
$adidas.Cells.Item(5, 3).Value = 'Promotional exclusivity'
$adidas.Cells.Item(5, 4).Value = 'Several ads offer an incentive to sign up for adiClub to participate in exclusive offers like winning a gift card. | This creates urgency around registration and exclusivity.'
$adidas.Cells.Item(5, 5).ClearContents()
$adidas.Cells.Item(5, 6).Value = 'Suscríbete a adiClub y crea tu lista de deseos para poder ganar una tarjeta regalo por valor de 250€ este Black Friday. (Sign up for adiClub and create your wishlist to win a €250 gift card this Black Friday.)'

$adidas.Cells.Item(6, 1).Value = 3
$adidas.Cells.Item(6, 2).Value = 'adidas Vibes Fragrance'
$adidas.Cells.Item(6, 3).Value = 'Category extension'
$adidas.Cells.Item(6, 4).Value = 'This ad highlights adidas''s entry into the fragrance market, positioning it as part of the lifestyle beyond just apparel. | This suggests adidas is becoming synonymous with an all-around lifestyle brand.'
$adidas.Cells.Item(6, 5).ClearContents()
$adidas.Cells.Item(6, 6).Value = 'adidas Vibes Body & Hair Mist mein Frische-Boost für jeden Mood, ob vor dem Training oder beim Entspannen zuhause #adidasfragrance #createdwithadidas (adidas Vibes Body & Hair Mist, my freshness boost for any mood, whether before training or relaxing at home #adidasfragrance #createdwithadidas)'

# ===================================================================
# View-state touches matching the saved workbook: adidas tab becomes
# the active/selected tab (PUMA tab loses the selection).
# ===================================================================
$adidas.Activate()
$adidas.Range("I12").Select()